$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("I9")
$r.Interior.Pattern = 1
$r.Interior.Pattern = -4142
